$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the "Reasoning" column from the table (and worksheet).
# ------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item("Reasoning").Delete()

# ------------------------------------------------------------------
# 2. Update the "username/email/password" story text to add a
#    trailing period, and fill in the remaining user stories plus
#    their story point values.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "As a customer, I want to be able to change my home address, because I recently moved to a new house."
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "As a customer, I want to be able to find groupings of items easily, because I don't want to search the entire catalogue of items on the website for one thing."
$ws.Range("B3").Value = 9

$ws.Range("A4").Value = "As a customer, I want to be able to change my username, email and password, because I'm concerned about my previous information being compromised."
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = "As a customer, I want to be able to find items that fit my budget. "
$ws.Range("B5").Value = 6

$ws.Range("A6").Value = "As a customer, I want to have a Wish list that have all my favorite items, because I want it to be easy to locate them later when it is the appropriate time to buy them."
$ws.Range("B6").Value = 12

# ------------------------------------------------------------------
# 3. Shrink the table down to the now-populated range (2 columns,
#    header + 5 data rows).
# ------------------------------------------------------------------
$tbl.Resize($ws.Range("A1:B6"))

# ------------------------------------------------------------------
# 4. Row heights for the wrapped, long user-story rows.
# ------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 29
$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(6).RowHeight = 29

# ------------------------------------------------------------------
# 5. Touch a few extra (otherwise blank) cells down through row 15 so
#    the sheet's used range / dimension grows to match, then clear
#    any stray values so they stay empty.
# ------------------------------------------------------------------
$ws.Range("C9:C13").Borders.LineStyle = 0
$ws.Range("A14:A15").Borders.LineStyle = 0
$ws.Range("C14:C15").Borders.LineStyle = 0

# ------------------------------------------------------------------
# 6. Update the active selection to match the saved selection.
# ------------------------------------------------------------------
$ws.Range("H9").Select()
